$wb = $excel.ActiveWorkbook

# --- Create the new sheet ("Hoja2") and drop the old one ("Hoja1") ---------
# Doing it this way (Add + Delete) mirrors what real Excel does when a sheet
# is recreated, bumping the internal sheetId counter (1 -> 2), matching the
# target workbook.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Hoja2"
$oldSheet = $wb.Worksheets.Item("Hoja1")
$oldSheet.Delete()

$ws = $newSheet

# --- Header row --------------------------------------------------------
$ws.Range("A1").Value = "rut"
$ws.Range("B1").Value = "fecha"
$ws.Range("C1").Value = "puntaje"
$ws.Range("A1:C1").NumberFormat = "General"

# --- Data rows -----------------------------------------------------------
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "20.999.124-8"
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = (Get-Date -Year 2023 -Month 10 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C2").NumberFormat = "0.00"
$ws.Range("C2").Value = 950

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "rut1"
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = (Get-Date -Year 2023 -Month 10 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C3").NumberFormat = "0.00"
$ws.Range("C3").Value = 980

# --- Column width (best-fit on column B, like Excel's AutoFit) -----------
$ws.Columns.Item(2).ColumnWidth = 23.28515625

# --- Selection -------------------------------------------------------------
$ws.Range("E10").Select()
